$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 93
$ws.Range("H93").Value = 34159.3
$ws.Range("J93").Value = 34159.3
$ws.Range("L93").Value = 34159.3
$ws.Range("N93").Value = -39151.3
# Row 98
$ws.Range("H98").Value = 36387.5
$ws.Range("I98").Value = 987.4
$ws.Range("J98").Value = 390388.5
$ws.Range("K98").Value = 987.4
$ws.Range("L98").Value = 390388.5
$ws.Range("M98").Value = 510.6
$ws.Range("N98").Value = -393384.5
# Row 107
$ws.Range("H107").Value = 755.4545000000001
$ws.Range("I107").Value = 755.4545000000001
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 755.4545000000001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1164.5455
$ws.Range("N107").ClearContents()
# Row 122
$ws.Range("H122").Value = 36387.5
$ws.Range("I122").Value = 987.4
$ws.Range("J122").Value = 390388.5
$ws.Range("K122").Value = 2962.2
$ws.Range("L122").Value = 1171165.5
$ws.Range("M122").Value = -512.1999999999998
$ws.Range("N122").Value = -1176065.5
# Row 129
$ws.Range("H129").Value = 302296.2
$ws.Range("I129").Value = 751295.9
$ws.Range("J129").Value = 2963.0833
$ws.Range("K129").Value = 2253887.7
$ws.Range("L129").Value = 8889.249899999999
$ws.Range("M129").Value = -2248887.7
$ws.Range("N129").Value = -18889.2499
# Row 135
$ws.Range("I135").Value = 900.7742
$ws.Range("K135").Value = 8106.967799999999
$ws.Range("M135").Value = -5571.967799999999
# Row 138
$ws.Range("H138").Value = 2390.3718
$ws.Range("I138").Value = 2211.8276
$ws.Range("J138").Value = 2496.0408
$ws.Range("K138").Value = 6635.4828
$ws.Range("L138").Value = 7488.1224
$ws.Range("M138").Value = -1495.4828
$ws.Range("N138").Value = -17768.1224

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 38008.38
$ws.Range("I32").Value = 39304.66
$ws.Range("J32").Value = 28502.334
$ws.Range("K32").Value = 39304.66
$ws.Range("L32").Value = 28502.334
$ws.Range("M32").Value = -39017.66
$ws.Range("N32").Value = -29076.334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 131
$ws.Range("H131").Value = 25000
$ws.Range("J131").Value = 25000
$ws.Range("L131").Value = 25000
$ws.Range("N131").Value = -35080
# Row 134
$ws.Range("H134").Value = 1819.5435
$ws.Range("I134").Value = 1438.5385
$ws.Range("J134").Value = 3942.2856
$ws.Range("K134").Value = 4315.6155
$ws.Range("L134").Value = 11826.8568
$ws.Range("M134").Value = -1780.6155
$ws.Range("N134").Value = -16896.8568

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 35
$ws.Range("H35").Value = 2652.4
$ws.Range("I35").Value = 1440.5
$ws.Range("J35").Value = 7500
$ws.Range("K35").Value = 1440.5
$ws.Range("L35").Value = 7500
$ws.Range("M35").Value = -1146.5
$ws.Range("N35").Value = -8088
# Row 99
$ws.Range("H99").Value = 2156.7407
$ws.Range("I99").Value = 2122.4614
$ws.Range("J99").Value = 2188.5715
$ws.Range("K99").Value = 2122.4614
$ws.Range("L99").Value = 2188.5715
$ws.Range("M99").Value = -624.4614000000001
$ws.Range("N99").Value = -5184.5715
# Row 126
$ws.Range("H126").Value = 2156.7407
$ws.Range("I126").Value = 2122.4614
$ws.Range("J126").Value = 2188.5715
$ws.Range("K126").Value = 6367.3842
$ws.Range("L126").Value = 6565.7145
$ws.Range("M126").Value = -3897.3842
$ws.Range("N126").Value = -11505.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 4988.5654
$ws.Range("I113").Value = 8987.666999999999
$ws.Range("J113").Value = 625.9091
$ws.Range("K113").Value = 26963.001
$ws.Range("L113").Value = 1877.7273
$ws.Range("M113").Value = -24793.001
$ws.Range("N113").Value = -6217.7273
# Row 132
$ws.Range("H132").Value = 1850.625
$ws.Range("I132").Value = 1250
$ws.Range("J132").Value = 2451.25
$ws.Range("K132").Value = 11250
$ws.Range("L132").Value = 22061.25
$ws.Range("M132").Value = -8720
$ws.Range("N132").Value = -27121.25
# Row 137
$ws.Range("H137").Value = 41675656
$ws.Range("J137").Value = 90925030
$ws.Range("L137").Value = 272775090
$ws.Range("N137").Value = -272785290

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 691.1818
$ws.Range("I55").Value = 717.7273
$ws.Range("J55").Value = 664.63635
$ws.Range("K55").Value = 717.7273
$ws.Range("L55").Value = 664.63635
$ws.Range("M55").Value = -544.7273
$ws.Range("N55").Value = -1010.63635
# Row 61
$ws.Range("H61").Value = 2019.9354
$ws.Range("I61").Value = 2068.5
$ws.Range("K61").Value = 2068.5
$ws.Range("M61").Value = -1866.5
# Row 92
$ws.Range("H92").Value = 44389
$ws.Range("J92").Value = 44389
$ws.Range("L92").Value = 44389
$ws.Range("N92").Value = -49381
# Row 96
$ws.Range("H96").Value = 30037.8
$ws.Range("J96").Value = 30037.8
$ws.Range("L96").Value = 30037.8
$ws.Range("N96").Value = -35529.8
# Row 113
$ws.Range("H113").Value = 2019.9354
$ws.Range("I113").Value = 2068.5
$ws.Range("K113").Value = 2068.5
$ws.Range("M113").Value = 101.5
# Row 133
$ws.Range("H133").Value = 24836
$ws.Range("J133").Value = 24836
$ws.Range("L133").Value = 24836
$ws.Range("N133").Value = -29896
# Row 136
$ws.Range("H136").Value = 3279.8333
$ws.Range("I136").Value = 2914
$ws.Range("K136").Value = 8742
$ws.Range("M136").Value = -6192

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 62
$ws.Range("H62").Value = 2997.6667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2997.6667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2997.6667
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4245.6667
# Row 65
$ws.Range("H65").Value = 2997.6667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2997.6667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 14988.3335
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -21228.3335
# Row 126
$ws.Range("H126").Value = 2262981.8
$ws.Range("I126").Value = 2674251.2
$ws.Range("K126").Value = 8022753.600000001
$ws.Range("M126").Value = -8020283.600000001
# Row 135
$ws.Range("H135").Value = 52266.332
$ws.Range("J135").Value = 52266.332
$ws.Range("L135").Value = 52266.332
$ws.Range("N135").Value = -62406.332
# Row 136
$ws.Range("H136").Value = 21553.02
$ws.Range("I136").Value = 92687
$ws.Range("J136").Value = 2468.2927
$ws.Range("K136").Value = 278061
$ws.Range("L136").Value = 7404.8781
$ws.Range("M136").Value = -275511
$ws.Range("N136").Value = -12504.8781
# Row 137
$ws.Range("H137").Value = 62999.5
$ws.Range("J137").Value = 62999.5
$ws.Range("L137").Value = 62999.5
$ws.Range("N137").Value = -73199.5
